$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4 looks like a date ("06/08/2025") and Excel would normally auto-convert it
# to a date serial number on assignment. Force it to stay plain text by
# temporarily marking the cell as Text before writing the value, then clear
# the formatting override so the cell keeps the sheet's default style.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "06/08/2025"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = "Basel"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Young Boys"
$ws.Range("F4").Value = "L"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 0.47
$ws.Range("L4").Value = 2.81
$ws.Range("M4").Value = 13
$ws.Range("N4").Value = 22
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = 12
